# Apply "Penalty Reward System" edit to B082DZ8HLT_po_data.xlsx
# Sheet1 "Weekly Quantity": keep header + row2, replace row3 with 45137.99999999999 / 1,
#   delete old rows 4 and 5.
# Sheet2 "Monthly Trend": keep header + row2, replace row3 with 45138.99999999999 / 1,
#   delete old row 4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws1.Range("A3").Value = 45137.99999999999
$ws1.Range("B3").Value = 1
$ws1.Rows.Item(4).EntireRow.Delete()
$ws1.Rows.Item(4).EntireRow.Delete()

$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Range("A3").Value = 45138.99999999999
$ws2.Range("B3").Value = 1
$ws2.Rows.Item(4).EntireRow.Delete()
